# "added a v2 for dutch grid"
# - CBC sheet: delivery_end (col D) rows 2-3 bumped 23 -> 24; becomes the active sheet; selection -> D5
# - RD sheet: delivery_end (col D) rows 2-7 bumped 23 -> 24; delivery_start (col C) rows 2-7 1 -> 0;
#             power (col G) row2 200 -> 20000, row5 500 -> 50000; no longer the active sheet; selection -> E13
# - lines sheet: selection -> G10

$wb = $excel.ActiveWorkbook

# --- "lines" sheet: just a selection move, no data change ---
$wsLines = $wb.Worksheets.Item("lines")
$wsLines.Range("G10").Select() | Out-Null

# --- "RD" sheet: bump delivery window + two power values, move selection ---
$wsRD = $wb.Worksheets.Item("RD")
$wsRD.Range("C2:C7").Value = 0
$wsRD.Range("D2:D7").Value = 24
$wsRD.Range("G2").Value = 20000
$wsRD.Range("G5").Value = 50000
$wsRD.Range("E13").Select() | Out-Null

# --- "CBC" sheet: new v2 delivery window, becomes the active tab ---
$wsCBC = $wb.Worksheets.Item("CBC")
$wsCBC.Range("D2:D3").Value = 24
$wsCBC.Activate() | Out-Null
$wsCBC.Range("D5").Select() | Out-Null
